$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Delete trailing rows 74-78 (now removed from the sheet) ---
$ws.Range("A74:A78").EntireRow.Delete()

# --- Row 1: add trailing sequence cells CX1,CY1 ---
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 3
$ws.Range("D1").Value = 4
$ws.Range("E1").Value = 5
$ws.Range("F1").Value = 6
$ws.Range("G1").Value = 7
$ws.Range("H1").Value = 8
$ws.Range("I1").Value = 9
$ws.Range("J1").Value = 10
$ws.Range("K1").Value = 11
$ws.Range("L1").Value = 12
$ws.Range("M1").Value = 13
$ws.Range("N1").Value = 14
$ws.Range("O1").Value = 15
$ws.Range("P1").Value = 16
$ws.Range("Q1").Value = 17
$ws.Range("R1").Value = 18
$ws.Range("S1").Value = 19
$ws.Range("T1").Value = 20
$ws.Range("U1").Value = 21
$ws.Range("V1").Value = 22
$ws.Range("W1").Value = 23
$ws.Range("X1").Value = 24
$ws.Range("Y1").Value = 25
$ws.Range("Z1").Value = 26
$ws.Range("AA1").Value = 27
$ws.Range("AB1").Value = 28
$ws.Range("AC1").Value = 29
$ws.Range("AD1").Value = 30
$ws.Range("AE1").Value = 31
$ws.Range("AF1").Value = 32
$ws.Range("AG1").Value = 33
$ws.Range("AH1").Value = 34
$ws.Range("AI1").Value = 35
$ws.Range("AJ1").Value = 36
$ws.Range("AK1").Value = 37
$ws.Range("AL1").Value = 38
$ws.Range("AM1").Value = 39
$ws.Range("AN1").Value = 40
$ws.Range("AO1").Value = 41
$ws.Range("AP1").Value = 42
$ws.Range("AQ1").Value = 43
$ws.Range("AR1").Value = 44
$ws.Range("AS1").Value = 45
$ws.Range("AT1").Value = 46
$ws.Range("AU1").Value = 47
$ws.Range("AV1").Value = 48
$ws.Range("AW1").Value = 49
$ws.Range("AX1").Value = 50
$ws.Range("AY1").Value = 51
$ws.Range("AZ1").Value = 52
$ws.Range("BA1").Value = 53
$ws.Range("BB1").Value = 54
$ws.Range("BC1").Value = 55
$ws.Range("BD1").Value = 56
$ws.Range("BE1").Value = 57
$ws.Range("BF1").Value = 58
$ws.Range("BG1").Value = 59
$ws.Range("BH1").Value = 60
$ws.Range("BI1").Value = 61
$ws.Range("BJ1").Value = 62
$ws.Range("BK1").Value = 63
$ws.Range("BL1").Value = 64
$ws.Range("BM1").Value = 65
$ws.Range("BN1").Value = 66
$ws.Range("BO1").Value = 67
$ws.Range("BP1").Value = 68
$ws.Range("BQ1").Value = 69
$ws.Range("BR1").Value = 70
$ws.Range("BS1").Value = 71
$ws.Range("BT1").Value = 72
$ws.Range("BU1").Value = 73
$ws.Range("BV1").Value = 74
$ws.Range("BW1").Value = 75
$ws.Range("BX1").Value = 76
$ws.Range("BY1").Value = 77
$ws.Range("BZ1").Value = 78
$ws.Range("CA1").Value = 79
$ws.Range("CB1").Value = 80
$ws.Range("CC1").Value = 81
$ws.Range("CD1").Value = 82
$ws.Range("CE1").Value = 83
$ws.Range("CF1").Value = 84
$ws.Range("CG1").Value = 85
$ws.Range("CH1").Value = 86
$ws.Range("CI1").Value = 87
$ws.Range("CJ1").Value = 88
$ws.Range("CK1").Value = 89
$ws.Range("CL1").Value = 90
$ws.Range("CM1").Value = 91
$ws.Range("CN1").Value = 92
$ws.Range("CO1").Value = 93
$ws.Range("CP1").Value = 94
$ws.Range("CQ1").Value = 95
$ws.Range("CR1").Value = 96
$ws.Range("CS1").Value = 97
$ws.Range("CT1").Value = 98
$ws.Range("CU1").Value = 99
$ws.Range("CV1").Value = 100
$ws.Range("CW1").Value = 101
$ws.Range("CX1").Value = 102
$ws.Range("CY1").Value = 103

# --- Row 2: header labels shifted by 2 inserted columns (Timezone @ G, z @ J) ---
$ws.Range("A2").Value = "Grid"
$ws.Range("B2").Value = "Year"
$ws.Range("C2").Value = "StartDLS"
$ws.Range("D2").Value = "EndDLS"
$ws.Range("E2").Value = "lat"
$ws.Range("F2").Value = "lng"
$ws.Range("G2").Value = "Timezone"
$ws.Range("H2").Value = "SurfaceArea"
$ws.Range("I2").Value = "Alt"
$ws.Range("J2").Value = "z"
$ws.Range("K2").Value = "id"
$ws.Range("L2").Value = "ih"
$ws.Range("M2").Value = "imin"
$ws.Range("N2").Value = "Fr_Paved"
$ws.Range("O2").Value = "Fr_Bldgs"
$ws.Range("P2").Value = "Fr_EveTr"
$ws.Range("Q2").Value = "Fr_DecTr"
$ws.Range("R2").Value = "Fr_Grass"
$ws.Range("S2").Value = "Fr_Bsoil"
$ws.Range("T2").Value = "Fr_Water"
$ws.Range("U2").Value = "IrrFr_EveTr"
$ws.Range("V2").Value = "IrrFr_DecTr"
$ws.Range("W2").Value = "IrrFr_Grass"
$ws.Range("X2").Value = "H_Bldgs"
$ws.Range("Y2").Value = "H_EveTr"
$ws.Range("Z2").Value = "H_DecTr"
$ws.Range("AA2").Value = "z0"
$ws.Range("AB2").Value = "zd"
$ws.Range("AC2").Value = "FAI_Bldgs"
$ws.Range("AD2").Value = "FAI_EveTr"
$ws.Range("AE2").Value = "FAI_DecTr"
$ws.Range("AF2").Value = "PopDensDay"
$ws.Range("AG2").Value = "PopDensNight"
$ws.Range("AH2").Value = "TrafficRate"
$ws.Range("AI2").Value = "BuildEnergyUse"
$ws.Range("AJ2").Value = "Code_Paved"
$ws.Range("AK2").Value = "Code_Bldgs"
$ws.Range("AL2").Value = "Code_EveTr"
$ws.Range("AM2").Value = "Code_DecTr"
$ws.Range("AN2").Value = "Code_Grass"
$ws.Range("AO2").Value = "Code_Bsoil"
$ws.Range("AP2").Value = "Code_Water"
$ws.Range("AQ2").Value = "LUMPS_DrRate"
$ws.Range("AR2").Value = "LUMPS_Cover"
$ws.Range("AS2").Value = "LUMPS_MaxRes"
$ws.Range("AT2").Value = "NARP_Trans"
$ws.Range("AU2").Value = "CondCode"
$ws.Range("AV2").Value = "SnowCode"
$ws.Range("AW2").Value = "SnowClearingProfWD"
$ws.Range("AX2").Value = "SnowClearingProfWE"
$ws.Range("AY2").Value = "AnthropogenicCode"
$ws.Range("AZ2").Value = "EnergyUseProfWD"
$ws.Range("BA2").Value = "EnergyUseProfWE"
$ws.Range("BB2").Value = "ActivityProfWD"
$ws.Range("BC2").Value = "ActivityProfWE"
$ws.Range("BD2").Value = "IrrigationCode"
$ws.Range("BE2").Value = "WaterUseProfManuWD"
$ws.Range("BF2").Value = "WaterUseProfManuWE"
$ws.Range("BG2").Value = "WaterUseProfAutoWD"
$ws.Range("BH2").Value = "WaterUseProfAutoWE"
$ws.Range("BI2").Value = "FlowChange"
$ws.Range("BJ2").Value = "RunoffToWater"
$ws.Range("BK2").Value = "PipeCapacity"
$ws.Range("BL2").Value = "GridConnection1of8"
$ws.Range("BM2").Value = "Fraction1of8"
$ws.Range("BN2").Value = "GridConnection2of8"
$ws.Range("BO2").Value = "Fraction2of8"
$ws.Range("BP2").Value = "GridConnection3of8"
$ws.Range("BQ2").Value = "Fraction3of8"
$ws.Range("BR2").Value = "GridConnection4of8"
$ws.Range("BS2").Value = "Fraction4of8"
$ws.Range("BT2").Value = "GridConnection5of8"
$ws.Range("BU2").Value = "Fraction5of8"
$ws.Range("BV2").Value = "GridConnection6of8"
$ws.Range("BW2").Value = "Fraction6of8"
$ws.Range("BX2").Value = "GridConnection7of8"
$ws.Range("BY2").Value = "Fraction7of8"
$ws.Range("BZ2").Value = "GridConnection8of8"
$ws.Range("CA2").Value = "Fraction8of8"
$ws.Range("CB2").Value = "WithinGridPavedCode"
$ws.Range("CC2").Value = "WithinGridBldgsCode"
$ws.Range("CD2").Value = "WithinGridEveTrCode"
$ws.Range("CE2").Value = "WithinGridDecTrCode"
$ws.Range("CF2").Value = "WithinGridGrassCode"
$ws.Range("CG2").Value = "WithinGridUnmanBSoilCode"
$ws.Range("CH2").Value = "WithinGridWaterCode"
$ws.Range("CI2").Value = "AreaWall"
$ws.Range("CJ2").Value = "Fr_ESTMClass_Paved1"
$ws.Range("CK2").Value = "Fr_ESTMClass_Paved2"
$ws.Range("CL2").Value = "Fr_ESTMClass_Paved3"
$ws.Range("CM2").Value = "Code_ESTMClass_Paved1"
$ws.Range("CN2").Value = "Code_ESTMClass_Paved2"
$ws.Range("CO2").Value = "Code_ESTMClass_Paved3"
$ws.Range("CP2").Value = "Fr_ESTMClass_Bldgs1"
$ws.Range("CQ2").Value = "Fr_ESTMClass_Bldgs2"
$ws.Range("CR2").Value = "Fr_ESTMClass_Bldgs3"
$ws.Range("CS2").Value = "Fr_ESTMClass_Bldgs4"
$ws.Range("CT2").Value = "Fr_ESTMClass_Bldgs5"
$ws.Range("CU2").Value = "Code_ESTMClass_Bldgs1"
$ws.Range("CV2").Value = "Code_ESTMClass_Bldgs2"
$ws.Range("CW2").Value = "Code_ESTMClass_Bldgs3"
$ws.Range("CX2").Value = "Code_ESTMClass_Bldgs4"
$ws.Range("CY2").Value = "Code_ESTMClass_Bldgs5"
$ws.Range("CZ2").Value = "!"
$ws.Range("DA2").Value = "Site"
$ws.Range("DB2").Value = "Reference"

# --- Row 3: data values updated + shifted + 3 new trailing cells ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2012
$ws.Range("C3").Value = 85
$ws.Range("D3").Value = 302
$ws.Range("E3").Value = 51.51
$ws.Range("F3").Value = 0.12
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 2221.56
$ws.Range("I3").Value = 10.7
$ws.Range("J3").Value = 49.6
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0.43
$ws.Range("O3").Value = 0.38
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0.02
$ws.Range("R3").Value = 0.03
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0.14000000000000001
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 22
$ws.Range("Y3").Value = 13.1
$ws.Range("Z3").Value = 13.1
$ws.Range("AA3").Value = 1.9
$ws.Range("AB3").Value = 14.2
$ws.Range("AC3").Value = -999
$ws.Range("AD3").Value = -999
$ws.Range("AE3").Value = -999
$ws.Range("AF3").Value = -999
$ws.Range("AG3").Value = 204.58
$ws.Range("AH3").Value = -999
$ws.Range("AI3").Value = -999
$ws.Range("AJ3").Value = 661
$ws.Range("AK3").Value = 662
$ws.Range("AL3").Value = 661
$ws.Range("AM3").Value = 662
$ws.Range("AN3").Value = 663
$ws.Range("AO3").Value = 663
$ws.Range("AP3").Value = 661
$ws.Range("AQ3").Value = 0.25
$ws.Range("AR3").Value = 1
$ws.Range("AS3").Value = 10
$ws.Range("AT3").Value = 1
$ws.Range("AU3").Value = 200
$ws.Range("AV3").Value = 99999
$ws.Range("AW3").Value = 99999
$ws.Range("AX3").Value = 99999
$ws.Range("AY3").Value = 661
$ws.Range("AZ3").Value = 661
$ws.Range("BA3").Value = 662
$ws.Range("BB3").Value = 90000
$ws.Range("BC3").Value = 90000
$ws.Range("BD3").Value = 99999
$ws.Range("BE3").Value = 99999
$ws.Range("BF3").Value = 99999
$ws.Range("BG3").Value = 99999
$ws.Range("BH3").Value = 99999
$ws.Range("BI3").Value = 0
$ws.Range("BJ3").Value = 0
$ws.Range("BK3").Value = 100
$ws.Range("BL3").Value = 0
$ws.Range("BM3").Value = 0
$ws.Range("BN3").Value = 0
$ws.Range("BO3").Value = 0
$ws.Range("BP3").Value = 0
$ws.Range("BQ3").Value = 0
$ws.Range("BR3").Value = 0
$ws.Range("BS3").Value = 0
$ws.Range("BT3").Value = 0
$ws.Range("BU3").Value = 0
$ws.Range("BV3").Value = 0
$ws.Range("BW3").Value = 0
$ws.Range("BX3").Value = 0
$ws.Range("BY3").Value = 0
$ws.Range("BZ3").Value = 0
$ws.Range("CA3").Value = 0
$ws.Range("CB3").Value = 661
$ws.Range("CC3").Value = 662
$ws.Range("CD3").Value = 663
$ws.Range("CE3").Value = 664
$ws.Range("CF3").Value = 665
$ws.Range("CG3").Value = 666
$ws.Range("CH3").Value = 667
$ws.Range("CI3").Value = 1.08
$ws.Range("CJ3").Value = 0.15
$ws.Range("CK3").Value = 0.05
$ws.Range("CL3").Value = 0.8
$ws.Range("CM3").Value = 806
$ws.Range("CN3").Value = 807
$ws.Range("CO3").Value = 808
$ws.Range("CP3").Value = 0.15
$ws.Range("CQ3").Value = 0.45
$ws.Range("CR3").Value = 0.05
$ws.Range("CS3").Value = 0.1
$ws.Range("CT3").Value = 0.25
$ws.Range("CU3").Value = 801
$ws.Range("CV3").Value = 802
$ws.Range("CW3").Value = 803
$ws.Range("CX3").Value = 804
$ws.Range("CY3").Value = 805
$ws.Range("CZ3").Value = "!"
$ws.Range("DA3").Value = "London"
$ws.Range("DB3").Value = "Kotthaus and Grimmond (2013, 2014a, 2014b)"

# --- Selection matches target sheetView ---
$ws.Range("B11").Select()

